{"js": "// Update the date line and all 25 division-problem answers in the table.\n// The new values replace the old ones in strict document order; some new\n// values happen to be identical to *other* (not-yet-edited) old values\n// elsewhere in the document, so we must not re-search the whole document\n// by text after each edit (that could match the wrong, already-updated\n// cell). Instead we resolve every target paragraph's index up front and\n// then write the new text into each one by that fixed index.\nconst newTexts = [\n  \"2025-08-23 Saturday\",\n  \"32\u00f76=5, 2\",\n  \"99\u00f72=49, 1\",\n  \"88\u00f76=14, 4\",\n  \"84\u00f77=12, 0\",\n  \"35\u00f76=5, 5\",\n  \"69\u00f77=9, 6\",\n  \"77\u00f74=19, 1\",\n  \"65\u00f74=16, 1\",\n  \"33\u00f78=4, 1\",\n  \"91\u00f77=13, 0\",\n  \"32\u00f73=10, 2\",\n  \"78\u00f75=15, 3\",\n  \"40\u00f76=6, 4\",\n  \"51\u00f78=6, 3\",\n  \"27\u00f73=9, 0\",\n  \"71\u00f74=17, 3\",\n  \"44\u00f76=7, 2\",\n  \"24\u00f73=8, 0\",\n  \"75\u00f76=12, 3\",\n  \"29\u00f75=5, 4\",\n  \"89\u00f74=22, 1\",\n  \"45\u00f72=22, 1\",\n  \"25\u00f72=12, 1\",\n  \"60\u00f73=20, 0\",\n  \"12\u00f74=3, 0\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Collect, in document order, the paragraphs that currently hold text\n// (the date line plus the 25 non-blank table cells). That is exactly the\n// same order the diff's old values appear in, so the Nth non-blank\n// paragraph maps to newTexts[N].\nconst targets = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text && paragraphs.items[i].text.length > 0) {\n    targets.push(paragraphs.items[i]);\n  }\n}\n\nif (targets.length !== newTexts.length) {\n  throw new Error(\n    \"Expected \" + newTexts.length + \" non-empty paragraphs, found \" + targets.length\n  );\n}\n\nfor (let i = 0; i < targets.length; i++) {\n  targets[i].insertText(newTexts[i], Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the date line and all 25 division-problem answers in the table.\n# Some new values coincide with other (not-yet-edited) old values elsewhere\n# in the document, so we resolve every target paragraph's index up front\n# (while the document still holds the original text) and only then write\n# the new text into each one by that fixed index. This avoids accidentally\n# re-matching an already-updated cell via a text search performed later.\n$d = $word.ActiveDocument\n\n$newTexts = @(\n    \"2025-08-23 Saturday\",\n    \"32\u00f76=5, 2\",\n    \"99\u00f72=49, 1\",\n    \"88\u00f76=14, 4\",\n    \"84\u00f77=12, 0\",\n    \"35\u00f76=5, 5\",\n    \"69\u00f77=9, 6\",\n    \"77\u00f74=19, 1\",\n    \"65\u00f74=16, 1\",\n    \"33\u00f78=4, 1\",\n    \"91\u00f77=13, 0\",\n    \"32\u00f73=10, 2\",\n    \"78\u00f75=15, 3\",\n    \"40\u00f76=6, 4\",\n    \"51\u00f78=6, 3\",\n    \"27\u00f73=9, 0\",\n    \"71\u00f74=17, 3\",\n    \"44\u00f76=7, 2\",\n    \"24\u00f73=8, 0\",\n    \"75\u00f76=12, 3\",\n    \"29\u00f75=5, 4\",\n    \"89\u00f74=22, 1\",\n    \"45\u00f72=22, 1\",\n    \"25\u00f72=12, 1\",\n    \"60\u00f73=20, 0\",\n    \"12\u00f74=3, 0\"\n)\n\n# Collect, in document order, the paragraphs that currently hold text (the\n# date line plus the 25 non-blank table cells). That is exactly the same\n# order the diff's old values appear in, so the Nth non-blank paragraph\n# maps to $newTexts[N-1].\n$paragraphCount = $d.Paragraphs.Count\n$targets = @()\nfor ($i = 1; $i -le $paragraphCount; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $trimmed = $p.Range.Text.Trim([char]13, [char]7, [char]10)\n    if ($trimmed.Length -gt 0) {\n        $targets += $p\n    }\n}\n\nif ($targets.Count -ne $newTexts.Count) {\n    throw \"Expected $($newTexts.Count) non-empty paragraphs, found $($targets.Count)\"\n}\n\nfor ($i = 0; $i -lt $targets.Count; $i++) {\n    $targets[$i].Range.Text = $newTexts[$i]\n}\n"}
